$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking strings (e.g. "26.424.44", "0.0606") that must
# stay literal text, not be reinterpreted as numbers/dates. Format the whole
# data range of column D as Text before writing, then restore the original
# (default/unstyled) look so the saved style table matches the source workbook.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.424.44"
$ws.Range("E2").Value = "  +0.58%  "
$ws.Range("D3").Value = "1.611.79"
$ws.Range("E3").Value = "  +1.16%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "211.98"
$ws.Range("E5").Value = "  -0.55%  "
$ws.Range("D6").Value = "0.496"
$ws.Range("E6").Value = "  -0.52%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("B8").Value = "Dogecoin"
$ws.Range("C8").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D8").Value = "0.0606"
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").Value = "0.244"
$ws.Range("E9").Value = "  -0.42%  "
$ws.Range("E10").Value = "  +1.30%  "
$ws.Range("E11").Value = "  -0.54%  "
$ws.Range("D12").Value = "1.839.29"
$ws.Range("E12").Value = "  +1.20%  "
$ws.Range("D13").Value = "1.617.20"
$ws.Range("E13").Value = "  +1.48%  "
$ws.Range("E14").Value = "  -0.18%  "
$ws.Range("E15").Value = "  -0.06%  "
$ws.Range("D16").Value = "63.55"
$ws.Range("E16").Value = "  -0.37%  "
$ws.Range("D17").Value = "234.88"
$ws.Range("E17").Value = "  +8.87%  "
$ws.Range("D18").Value = "26.432.17"
$ws.Range("E19").Value = "  +0.35%  "
$ws.Range("E20").Value = "  +3.75%  "
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("E22").Value = "  -0.47%  "
$ws.Range("E23").Value = "  +4.63%  "
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("D25").Value = "146.61"
$ws.Range("E25").Value = "  +1.31%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("E27").Value = "  +0.36%  "
$ws.Range("E28").Value = "  +0.11%  "
$ws.Range("E30").Value = "  +1.16%  "
$ws.Range("E31").Value = "  -0.36%  "
$ws.Range("D32").Value = "1.502.49"
$ws.Range("E32").Value = "  +6.65%  "
$ws.Range("D33").Value = "3.23"
$ws.Range("E33").Value = "  +1.24%  "
$ws.Range("D34").Value = "2.95"
$ws.Range("E34").Value = "  -0.96%  "
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").Value = "1.50"
$ws.Range("E35").Value = "  +2.79%  "
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").Value = "2.42"
$ws.Range("E36").Value = "  -0.20%  "
$ws.Range("D37").Value = "0.558"
$ws.Range("E37").Value = "  -2.86%  "
$ws.Range("D38").Value = "0.0164"
$ws.Range("E38").Value = "  -0.11%  "
$ws.Range("E39").Value = "  +0.30%  "
$ws.Range("D40").Value = "5.81"
$ws.Range("E40").Value = "  +0.70%  "
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("E42").Value = "  +0.83%  "
$ws.Range("D43").Value = "1.751.45"
$ws.Range("E43").Value = "  +1.30%  "
$ws.Range("D44").Value = "0.925"
$ws.Range("E44").Value = "  -3.92%  "
$ws.Range("D45").Value = "0.760"
$ws.Range("E45").Value = "  -0.14%  "
$ws.Range("D46").Value = "61.67"
$ws.Range("E46").Value = "  +1.36%  "
$ws.Range("D47").Value = "89.68"
$ws.Range("E47").Value = "  +2.89%  "
$ws.Range("E48").Value = "  -0.06%  "
$ws.Range("E49").Value = "  -0.04%  "
$ws.Range("E50").Value = "  +0.98%  "
$ws.Range("D51").Value = "7.46"
$ws.Range("E51").Value = "  +1.17%  "

$ws.Range("D2:D51").Style = "Normal"
